$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that currently reads:
#   "Las descentralización de datos y el aislamiento de procesos."
# ---------------------------------------------------------------------------
$oldFull = "Las descentralización de datos y el aislamiento de procesos."

$search = $d.Content
$found = $search.Find.Execute($oldFull, $false, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target paragraph text."
}

$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $search.Start -and $p.Range.End -ge $search.End) {
        $para = $p
    }
}
if ($para -eq $null) {
    throw "Could not resolve the containing paragraph."
}

# ---------------------------------------------------------------------------
# Replace the whole paragraph (keeping its original paragraph-level rsid
# attributes/pPr) with the edited wording:
#   - "Las descentralización de datos" stays as its own run
#   - the (previously trailing) "_GoBack" bookmark now sits right after it
#   - the rest of the sentence is rewritten/expanded, with "agil" wrapped in
#     spell-check proofErr markers, matching Word's own as-you-type markup
# ---------------------------------------------------------------------------
$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005D3FCA" w:rsidRPr="002B5097" w:rsidRDefault="002B5097" w:rsidP="005D3FCA">' +
    '<w:pPr>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
        '<w:ind w:left="426"/>' +
        '<w:rPr><w:u w:val="single"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:t>Las descentralización de datos</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve">, integración entre sistemas legados, el aislamiento de procesos, tiempo </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>agil</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> para el cambio y facilidad de este.</w:t></w:r>' +
    '</w:p>'

$para.Range.InsertXML($newParagraphXml) | Out-Null

Write-Output "done"
